# Restored from revision of DEFAULT on 12/20/2021 04:05:07 PM. Type: SAVE.
#
# C7 currently holds the shared string "From"; it becomes a brand-new
# string "From111". D8's numeric value 113 becomes 1133.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "From111"
$ws.Range("D8").Value = 1133
